$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column width adjustments (cols D, E, I, K widened)
# Target stored widths (OOXML "width" attribute, character units) are achieved by
# compensating for the engine's pixel-snap rounding (snaps to nearest 1/7 char, +5px pad).
$ws.Columns.Item(4).ColumnWidth = 12.0
$ws.Columns.Item(5).ColumnWidth = 12.142857142857142
$ws.Columns.Item(9).ColumnWidth = 67.42857142857143
$ws.Columns.Item(11).ColumnWidth = 30.428571428571427

# Insert a thin space before "%" in English indicator descriptions (columns I and K),
# and drop a trailing period in I63.
$ws.Range("I5").Value = "Increase the proportion of organically farmed agricultural land to 30 % by 2030"
$ws.Range("K5").Value = "increase the proportion of organically farmed agricultural land to 30 % by 2030"
$ws.Range("I9").Value = "Reduction to 7 % by 2030"
$ws.Range("K9").Value = "reduction to 7 % by 2030"
$ws.Range("I10").Value = "Reduction to 19 % by 2030"
$ws.Range("K10").Value = "reduction to 19 % by 2030"
$ws.Range("I13").Value = "Reduction of emissions to 55 % of 2005 level (unweighted average of the five pollutants) by 2030"
$ws.Range("K13").Value = "reduction to 55 % by 2030"
$ws.Range("I16").Value = "Reduce the proportion to 9.5 % by 2030"
$ws.Range("K16").Value = "reduce the proportion to 9.5 % by 2030"
$ws.Range("I17").Value = "Increase the proportion to 55 % by 2030"
$ws.Range("K17").Value = "increase the proportion to 55 % by 2030"
$ws.Range("I18").Value = "Increase to 35 % by 2030"
$ws.Range("K18").Value = "increase to 35 % by 2030"
$ws.Range("I19").Value = "Increase to 70 % by 2030"
$ws.Range("K19").Value = "increase to 70 % by 2030"
$ws.Range("I20").Value = "Reduce the gap to 10 % by 2020, maintained until 2030 subsequently"
$ws.Range("K20").Value = "maintaine a maximum of 10 % by 2030"
$ws.Range("I21").Value = "30 % women on supervisory boards of listed and fully co-determined companies by 2030"
$ws.Range("I23").Value = "65 % by 2030"
$ws.Range("K23").Value = "65 % by 2030"
$ws.Range("I30").Value = "Increase by 2.1 % per year from 2008 to 2050"
$ws.Range("K30").Value = "increase by 2.1 % per year"
$ws.Range("I31").Value = "Reduction by 20 % by 2020, by 30 % by 2030, and by 50 % by 2050, all compared to 2008"
$ws.Range("K31").Value = "reduction by 30 % by 2030"
$ws.Range("I32").Value = "Increase to 18 % by 2020 and to 30 % by 2030, to 45 % by 2040 and to 60 % by 2050"
$ws.Range("K32").Value = "increase to 30 % by 2030"
$ws.Range("I33").Value = "Increase to at least 80 % by 2030"
$ws.Range("K33").Value = "increase to 80 % by 2030"
$ws.Range("I35").Value = "Annual government deficit less than 3 % of GDP, to be maintained until 2030"
$ws.Range("K35").Value = "less than 3 % of GDP"
$ws.Range("I36").Value = "Structurally balanced government budget, general government structural deficit must not exceed 0.5 % of GDP, to be maintained until 2030"
$ws.Range("K36").Value = "less than 0.5 % of GDP"
$ws.Range("I37").Value = "Ratio of government debt to GDP must not exceed  60 %, to be maintained until 2030"
$ws.Range("K37").Value = "max. 60 % of GDP"
$ws.Range("I40").Value = "Increase to 78 % by 2030"
$ws.Range("K40").Value = "increase to 78 % by 2030"
$ws.Range("I41").Value = "Increase to 60 % by 2030"
$ws.Range("K41").Value = "increase to 60 % by 2030"
$ws.Range("I43").Value = "At least 3.5 % of GDP per year by 2025"
$ws.Range("K43").Value = "at least 3.5 % of GDP per year by 2025"
$ws.Range("I50").Value = "Reduction by 15–20 % by 2030"
$ws.Range("K50").Value = "reduction by 15-20 % by 2030"
$ws.Range("I51").Value = "Reduction by 15–20 % by 2030"
$ws.Range("K51").Value = "reduction by 15-20 % by 2030"
$ws.Range("I53").Value = "Reduce the proportion of people who are overburdened to 13 % by 2030"
$ws.Range("K53").Value = "reduce to 13 % by 2030"
$ws.Range("I55").Value = "Increase the market share to 34 % by 2030"
$ws.Range("K55").Value = "increase the market share to 34 % by 2030"
$ws.Range("I60").Value = "Increase the proportion to 95 % by 2020"
$ws.Range("K60").Value = "increase to 95 % by 2020"
$ws.Range("I62").Value = "Reduce by at least 65 % by 2030 and by at least 88 % by 2040; greenhouse gas neutrality to be achieved by 2045"
$ws.Range("K62").Value = "reduce by 65 % by 2030"
$ws.Range("I63").Value = "Increase international climate finance to at least 6 billion euros by 2025 at the latest"
$ws.Range("I68").Value = "Reduction by 35 % by 2030 compared to 2005"
$ws.Range("K68").Value = "reduction by 35 % by 2030"
$ws.Range("I75").Value = "Increase the proportion to 0.7 % of gross national income by 2030"
$ws.Range("K75").Value = "increase to 0.7 % of GNI by 2030"
$ws.Range("I76").Value = "Increase the number by 10 % from 2015 to 2020, then stabilised"
$ws.Range("K76").Value = "increase by 10 % from 2015 to 2020, then stabilised"
